$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = 104806.3676698649
$ws.Range("E4").Value = -0.05427801192534044
$ws.Range("F4").Value = 0.2237048278216053
$ws.Range("G4").Value = -0.998740133840526
$ws.Range("H4").Value = 8.102398338075705
$ws.Range("D5").Value = 105442.7421606334
$ws.Range("E5").Value = -0.05695932407022027
$ws.Range("F5").Value = 0.2710939404004757
$ws.Range("G5").Value = -2.048719180784073
$ws.Range("H5").Value = 18.20988765032978
$ws.Range("D6").Value = 106042.563472038
$ws.Range("E6").Value = -0.06060798507591186
$ws.Range("F6").Value = 0.2747870632832989
$ws.Range("G6").Value = -1.542157623214067
$ws.Range("H6").Value = 12.22864035222789
$ws.Range("D7").Value = 106658.3775030086
$ws.Range("E7").Value = -0.0710663309355213
$ws.Range("F7").Value = 0.3330032235807257
$ws.Range("G7").Value = -2.238377426997002
$ws.Range("H7").Value = 15.60446483463334
$ws.Range("D8").Value = 107299.6524794298
$ws.Range("E8").Value = -0.06339223643104629
$ws.Range("F8").Value = 0.2745970157817861
$ws.Range("G8").Value = -0.947614428794478
$ws.Range("H8").Value = 7.276718644171352
$ws.Range("D9").Value = 108514.9872278934
$ws.Range("E9").Value = -0.08098588922748917
$ws.Range("F9").Value = 0.2615696624129245
$ws.Range("G9").Value = -1.156236635598374
$ws.Range("H9").Value = 7.78197903668239
$ws.Range("D10").Value = 110327.0895838107
$ws.Range("E10").Value = -0.1180677421501296
$ws.Range("F10").Value = 0.4031756866540213
$ws.Range("G10").Value = -1.604277559724129
$ws.Range("H10").Value = 9.180149728575055
$ws.Range("D11").Value = 111468.1096214049
$ws.Range("E11").Value = -0.1602817381780751
$ws.Range("F11").Value = 0.5088050387519121
$ws.Range("G11").Value = -1.648640794411797
$ws.Range("H11").Value = 7.501703896379907
$ws.Range("D12").Value = 112544.354143515
$ws.Range("E12").Value = -0.1717696007443023
$ws.Range("F12").Value = 0.5023922918757532
$ws.Range("G12").Value = -1.520841333433046
$ws.Range("H12").Value = 6.681400323468833
$ws.Range("D14").Value = 104098.9759193644
$ws.Range("E14").Value = -0.07644172311150114
$ws.Range("F14").Value = 0.1716495259302693
$ws.Range("G14").Value = -0.4637748782001339
$ws.Range("H14").Value = 7.425377512616324
$ws.Range("D16").Value = 103326.7708636718
$ws.Range("E16").Value = -0.1616264213003777
$ws.Range("F16").Value = 0.1721305297528183
$ws.Range("G16").Value = -0.9193210319563938
$ws.Range("H16").Value = 6.447565790854446
$ws.Range("D18").Value = 104101.039984516
$ws.Range("E18").Value = -0.06891200518171775
$ws.Range("F18").Value = 0.1678159015713327
$ws.Range("G18").Value = -0.4344410596768205
$ws.Range("H18").Value = 5.701134797319142
$ws.Range("D20").Value = 104211.616198819
$ws.Range("E20").Value = -0.0639012016028859
$ws.Range("F20").Value = 0.2316135982264532
$ws.Range("G20").Value = -2.922824638629686
$ws.Range("H20").Value = 35.52112564720923
$ws.Range("D21").Value = 104802.5009393485
$ws.Range("E21").Value = -0.05693102613717658
$ws.Range("F21").Value = 0.3052255104890294
$ws.Range("G21").Value = -3.410687159462758
$ws.Range("H21").Value = 29.25562509652999
